$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5789666666666667
$ws.Range("H2").Value = 1.7369
$ws.Range("I2").Value = 0.01523705650035473
$ws.Range("J2").Value = 0.01523705650035472
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.699658666666667
$ws.Range("N2").Value = 20.098976
$ws.Range("O2").Value = 0.1402150605386345
$ws.Range("P2").Value = 0.1402150605386345
$ws.Range("Q2").Value = 3.878879046044445
$ws.Range("R2").Value = 34.9099114144
$ws.Range("S2").Value = 0.002136464799627833
$ws.Range("T2").Value = 0.002136464799627832
$ws.Range("G3").Value = 0.5789666666666667
$ws.Range("H3").Value = 1.7369
$ws.Range("I3").Value = 0.01523705650035473
$ws.Range("J3").Value = 0.01523705650035472
$ws.Range("O3").Value = 0.6453289538613627
$ws.Range("P3").Value = 0.6453289538613627
$ws.Range("Q3").Value = 17.852240318
$ws.Range("R3").Value = 160.670162862
$ws.Range("S3").Value = 0.009832913731300392
$ws.Range("T3").Value = 0.00983291373130039
$ws.Range("G4").Value = 0.5789666666666667
$ws.Range("H4").Value = 1.7369
$ws.Range("I4").Value = 0.01523705650035473
$ws.Range("J4").Value = 0.01523705650035472
$ws.Range("O4").Value = 0.2144559856000028
$ws.Range("P4").Value = 0.2144559856000028
$ws.Range("Q4").Value = 5.932663906766667
$ws.Range("R4").Value = 53.3939751609
$ws.Range("S4").Value = 0.003267677969426503
$ws.Range("T4").Value = 0.003267677969426502
$ws.Range("I5").Value = 0.6545086962501954
$ws.Range("J5").Value = 0.6545086962501954
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.699658666666667
$ws.Range("N5").Value = 20.098976
$ws.Range("O5").Value = 0.1402150605386345
$ws.Range("P5").Value = 0.1402150605386345
$ws.Range("Q5").Value = 166.6174872607218
$ws.Range("R5").Value = 1499.557385346496
$ws.Range("S5").Value = 0.0917719764677839
$ws.Range("T5").Value = 0.09177197646778389
$ws.Range("I6").Value = 0.6545086962501954
$ws.Range("J6").Value = 0.6545086962501954
$ws.Range("O6").Value = 0.6453289538613627
$ws.Range("P6").Value = 0.6453289538613627
$ws.Range("S6").Value = 0.422373412244303
$ws.Range("T6").Value = 0.422373412244303
$ws.Range("I7").Value = 0.6545086962501954
$ws.Range("J7").Value = 0.6545086962501954
$ws.Range("O7").Value = 0.2144559856000028
$ws.Range("P7").Value = 0.2144559856000028
$ws.Range("S7").Value = 0.1403633075381085
$ws.Range("T7").Value = 0.1403633075381085
$ws.Range("I8").Value = 0.33025424724945
$ws.Range("J8").Value = 0.3302542472494499
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.699658666666667
$ws.Range("N8").Value = 20.098976
$ws.Range("O8").Value = 0.1402150605386345
$ws.Range("P8").Value = 0.1402150605386345
$ws.Range("Q8").Value = 84.0724243224569
$ws.Range("R8").Value = 756.651818902112
$ws.Range("S8").Value = 0.0463066192712228
$ws.Range("T8").Value = 0.04630661927122279
$ws.Range("I9").Value = 0.33025424724945
$ws.Range("J9").Value = 0.3302542472494499
$ws.Range("O9").Value = 0.6453289538613627
$ws.Range("P9").Value = 0.6453289538613627
$ws.Range("S9").Value = 0.2131226278857594
$ws.Range("T9").Value = 0.2131226278857593
$ws.Range("I10").Value = 0.33025424724945
$ws.Range("J10").Value = 0.3302542472494499
$ws.Range("O10").Value = 0.2144559856000028
$ws.Range("P10").Value = 0.2144559856000028
$ws.Range("Q10").Value = 128.5870044957563
$ws.Range("S10").Value = 0.07082500009246782
$ws.Range("T10").Value = 0.0708250000924678
